$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update input values (carbon study tidy)
$ws.Range("B3").Value = 15
$ws.Range("H3").Value = 15

# Add label and new "2*ks" computed cell
$ws.Range("D6").Value = "2.*ks"
$ws.Range("E6").Formula = "=2*E5"

# Restore selection to E5 as the active cell
$null = $ws.Range("E5").Select()
